$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the row heights that belong to each row position - they stay
# tied to the row number, not to the data that happens to occupy it.
$h48 = $ws.Rows(48).RowHeight
$h49 = $ws.Rows(49).RowHeight
$h50 = $ws.Rows(50).RowHeight
$h52 = $ws.Rows(52).RowHeight

# Remove the "محلول رينجر" (Ringer's solution) line item - this was row 48.
# Deleting the entire row shifts every following row up by one.
$ws.Range("A48").EntireRow.Delete()

# Restore the original per-position row heights (48/49/50 keep the
# heights they had before the delete; the old footer row's height now
# belongs to row 51).
$ws.Rows(48).RowHeight = $h48
$ws.Rows(49).RowHeight = $h49
$ws.Rows(50).RowHeight = $h50
$ws.Rows(51).RowHeight = $h52

# The report generator that produced the authoritative version does not
# perfectly re-flow every column after the deletion, so a handful of
# cells need to be corrected by hand to match the published output.
$ws.Range("A48").Value = 42
$ws.Range("C48").Value = "24.00"
$ws.Range("H48").Value = "محلول ملح"
$ws.Range("N48").Value = "28:0"

$ws.Range("A49").Value = 43

# Grand total line (previously row 51, now row 50) drops by the price of
# the removed item (120.00 -> total goes from 2042.565 to 1922.565).
$ws.Range("P50").Value = 1922.5650000000001

# Footer timestamp - new export time.
$ws.Range("A51").Value = "Thursday, 4 September, 2025 4:57 PM"

$wb.Save()
